# Auto-generated edit script: apply cell value updates per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 65
$ws.Range("F5").Value = 502
$ws.Range("F6").Value = 936
$ws.Range("F9").Value = 2225
$ws.Range("F10").Value = 637
$ws.Range("F11").Value = 300
$ws.Range("F13").Value = 1116
$ws.Range("F15").Value = 2248
$ws.Range("F16").Value = 694
$ws.Range("F17").Value = 13992
$ws.Range("F18").Value = 10
$ws.Range("F19").Value = 1305
$ws.Range("F22").Value = 141
$ws.Range("F24").Value = 151
$ws.Range("F25").Value = 99
$ws.Range("F26").Value = 43
$ws.Range("F29").Value = 7
$ws.Range("F31").Value = 33
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 6
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = 22
$ws.Range("F9").Value = 152
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 477
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 65
$ws.Range("F4").Value = 477
$ws.Range("F6").Value = 6
$ws.Range("F7").Value = 502
$ws.Range("F8").Value = 936
$ws.Range("F9").Value = 4
$ws.Range("F11").Value = 22
$ws.Range("F12").Value = 2225
$ws.Range("F13").Value = 637
$ws.Range("F14").Value = 300
$ws.Range("F18").Value = 1116
$ws.Range("F21").Value = 152
$ws.Range("F23").Value = 2248
$ws.Range("F24").Value = 694
$ws.Range("C25").Value = '广州·APH亚细亚同人ONLY'
$ws.Range("D25").Value = '鸿盛二路巨大创意产业园 巨大产业园·智汇港'
$ws.Range("E25").Value = '2024.10.03 10:00-10.03 16:30'
$ws.Range("F25").Value = 10
$ws.Range("G25").Value = 68
$ws.Range("H25").Value = 'https://show.bilibili.com/platform/detail.html?id=92322'
$ws.Range("I25").Value = '//i1.hdslb.com/bfs/openplatform/202409/D4UjGOw61725339026536.jpeg'
$ws.Range("C26").Value = '广州·平野宏周粉丝见面会'
$ws.Range("D26").Value = '金花街道中山七路333号1906科技圆区3号楼109-1铺、110-1铺、111-1铺 音乐唐人馆'
$ws.Range("E26").Value = '2024.10.03 11:30-10.03 17:00'
$ws.Range("F26").Value = 86
$ws.Range("G26").Value = 200
$ws.Range("H26").Value = 'https://show.bilibili.com/platform/detail.html?id=90504'
$ws.Range("I26").Value = '//i0.hdslb.com/bfs/openplatform/202408/xe9bWfYv1723095455030.jpeg'
$ws.Range("B27").Value = '2024-10-03'
$ws.Range("C27").Value = '广州·音爆ANISON「超次元ACG室内音乐节」Vol.7'
$ws.Range("D27").Value = '恩宁路265号三层四层自编01 MAO Livehouse广州（永庆坊店）'
$ws.Range("E27").Value = '2024.10.03 16:30-10.03 23:00'
$ws.Range("F27").Value = 65
$ws.Range("G27").Value = 89
$ws.Range("H27").Value = 'https://show.bilibili.com/platform/detail.html?id=91653'
$ws.Range("I27").Value = '//i2.hdslb.com/bfs/openplatform/202408/QvbLfaqV1724923999529.jpeg'
$ws.Range("C28").Value = '广州·南部动漫节'
$ws.Range("D28").Value = '东沙大道16号 广州健康方舟'
$ws.Range("E28").Value = '2024.10.05 10:00-10.06 17:00'
$ws.Range("F28").Value = 1305
$ws.Range("H28").Value = 'https://show.bilibili.com/platform/detail.html?id=90923'
$ws.Range("I28").Value = '//i0.hdslb.com/bfs/openplatform/202409/UKoa3flf1726049294866.jpeg'
$ws.Range("B29").Value = '2024-10-05'
$ws.Range("C29").Value = '广州·文豪野犬同人only2.0'
$ws.Range("D29").Value = '会江路巨大产业园5栋2楼 国际会议中心'
$ws.Range("E29").Value = '2024.10.05 10:00-10.05 17:00'
$ws.Range("F29").Value = 65
$ws.Range("G29").Value = 60
$ws.Range("H29").Value = 'https://show.bilibili.com/platform/detail.html?id=92232'
$ws.Range("I29").Value = '//i1.hdslb.com/bfs/openplatform/202409/lf9TiiUG1725879597895.jpeg'
$ws.Range("C30").Value = '广州·凹凸世界ONLY'
$ws.Range("D30").Value = '广汕二路149号 飞梦篮球公园(高塘石馆)'
$ws.Range("F30").Value = 568
$ws.Range("G30").Value = 78
$ws.Range("H30").Value = 'https://show.bilibili.com/platform/detail.html?id=89715'
$ws.Range("I30").Value = '//i1.hdslb.com/bfs/openplatform/202407/BnOHWZD91721638002542.jpeg'
$ws.Range("C31").Value = '广州·樱漫潮玩动漫游戏嘉年华'
$ws.Range("D31").Value = '奥体南路12号 优托邦(奥体旗舰店)'
$ws.Range("F31").Value = 141
$ws.Range("G31").Value = 39.9
$ws.Range("H31").Value = 'https://show.bilibili.com/platform/detail.html?id=91453'
$ws.Range("I31").Value = '//i1.hdslb.com/bfs/openplatform/202408/kcavrR9W1724055693105.jpeg'
$ws.Range("B32").Value = '2024-10-06'
$ws.Range("C32").Value = '广州·第五人格同人only4.0'
$ws.Range("D32").Value = '会江路巨大产业园5栋2楼 国际会议中心'
$ws.Range("E32").Value = '2024.10.06 10:00-10.06 17:00'
$ws.Range("F32").Value = 38
$ws.Range("G32").Value = 60
$ws.Range("H32").Value = 'https://show.bilibili.com/platform/detail.html?id=92168'
$ws.Range("I32").Value = '//i0.hdslb.com/bfs/openplatform/202409/twg9GYOJ1725616149211.jpeg'
$ws.Range("B33").Value = '2024-10-07'
$ws.Range("C33").Value = '广州·第十届萌物语动漫嘉年华'
$ws.Range("D33").Value = '奥体南路12号 优托邦(奥体旗舰店)'
$ws.Range("E33").Value = '2024.10.07 10:00-10.07 17:00'
$ws.Range("F33").Value = 151
$ws.Range("G33").Value = 39.9
$ws.Range("H33").Value = 'https://show.bilibili.com/platform/detail.html?id=91162'
$ws.Range("I33").Value = '//i1.hdslb.com/bfs/openplatform/202408/9m6CSSzf1723739254235.jpeg'
$ws.Range("B34").Value = '2024-10-13'
$ws.Range("C34").Value = '广州·BanG Dream 同人Only'
$ws.Range("D34").Value = '华观路1932号 智慧城广百广场'
$ws.Range("E34").Value = '2024.10.13 10:00-10.13 18:00'
$ws.Range("F34").Value = 99
$ws.Range("G34").Value = 58
$ws.Range("H34").Value = 'https://show.bilibili.com/platform/detail.html?id=92314'
$ws.Range("I34").Value = '//i0.hdslb.com/bfs/openplatform/202409/HHkN2uUe1726045216331.png'
$ws.Range("C35").Value = '广州·Luca Stricagnoli 2024《进化时间》指弹吉他音乐会'
$ws.Range("D35").Value = '恩宁路265号3层 MaoLivehouse(永庆坊店)'
$ws.Range("E35").Value = '2024.10.19 19:30-10.19 21:00'
$ws.Range("F35").Value = 9
$ws.Range("G35").Value = 220
$ws.Range("H35").Value = 'https://show.bilibili.com/platform/detail.html?id=91352'
$ws.Range("I35").Value = '//i1.hdslb.com/bfs/openplatform/202408/ArhAlkP41724743278046.jpeg'
$ws.Range("B36").Value = '2024-10-19'
$ws.Range("C36").Value = '广州·次元喵喵动漫嘉年华02'
$ws.Range("D36").Value = '东沙大道16号 广州健康方舟'
$ws.Range("E36").Value = '2024.10.19 10:00-10.19 18:00'
$ws.Range("F36").Value = 43
$ws.Range("G36").Value = 60
$ws.Range("H36").Value = 'https://show.bilibili.com/platform/detail.html?id=91566'
$ws.Range("I36").Value = '//i1.hdslb.com/bfs/openplatform/202408/VJ9w4T6W1724046324480.jpeg'
$ws.Range("B37").Value = '2024-10-20'
$ws.Range("C37").Value = '广州·ROOKiEZ is PUNK`D 「Reignite Youth （重燃青春）」2024 CHINA Tour '
$ws.Range("D37").Value = '南洲路154号侨建大厦2F SDlivehouse'
$ws.Range("E37").Value = '2024.10.20 20:00-10.20 21:30'
$ws.Range("F37").Value = 7
$ws.Range("G37").Value = 259
$ws.Range("H37").Value = 'https://show.bilibili.com/platform/detail.html?id=92075'
$ws.Range("I37").Value = '//i2.hdslb.com/bfs/openplatform/202409/kAxVF2Jw1725542237304.jpeg'
$ws.Range("B38").Value = '2024-10-25'
$ws.Range("C38").Value = '广州·新生代流媒体小天后野田爱实 2024 巡演'
$ws.Range("D38").Value = '南洲路158号2F SD Livehouse'
$ws.Range("E38").Value = '2024.10.25 20:00-10.25 22:00'
$ws.Range("F38").Value = 5
$ws.Range("G38").Value = 280
$ws.Range("H38").Value = 'https://show.bilibili.com/platform/detail.html?id=91823'
$ws.Range("I38").Value = '//i0.hdslb.com/bfs/openplatform/202409/oN7FyQ8v1725347758464.jpeg'
$ws.Range("B39").Value = '2024-10-26'
$ws.Range("C39").Value = '广州·wio jumponly4.0万圣狂欢节'
$ws.Range("D39").Value = '黄边三横路一街1号 设计殿堂'
$ws.Range("E39").Value = '2024.10.26 10:00-10.27 17:00'
$ws.Range("F39").Value = 282
$ws.Range("G39").Value = 69.90000000000001
$ws.Range("H39").Value = 'https://show.bilibili.com/platform/detail.html?id=89588'
$ws.Range("I39").Value = '//i0.hdslb.com/bfs/openplatform/202407/2kN5bTGE1721377069804.png'
$ws.Range("B40").Value = '2024-10-27'
$ws.Range("C40").Value = '广州·卡农·世界经典音乐之旅交响音乐会'
$ws.Range("D40").Value = '东风中路299号 广州中山纪念堂'
$ws.Range("E40").Value = '2024.10.27 19:30-10.27 21:00'
$ws.Range("F40").Value = 4
$ws.Range("G40").Value = 75
$ws.Range("H40").Value = 'https://show.bilibili.com/platform/detail.html?id=91040'
$ws.Range("I40").Value = '//i1.hdslb.com/bfs/openplatform/202408/WEqD8aj31724134831558.jpeg'
$ws.Range("B41").Value = '2024-11-08'
$ws.Range("C41").Value = '广州·「心随歌行」KOKIA 2024 中国巡演'
$ws.Range("D41").Value = '广州大道中1229号 广东艺术剧院'
$ws.Range("E41").Value = '2024.11.08 19:30-11.08 21:30'
$ws.Range("F41").Value = 35
$ws.Range("G41").Value = 880
$ws.Range("H41").Value = 'https://show.bilibili.com/platform/detail.html?id=90392'
$ws.Range("I41").Value = '//i0.hdslb.com/bfs/openplatform/202408/FDsbokRk1722914443578.jpeg'
$ws.Range("B42").Value = '2024-11-17'
$ws.Range("C42").Value = '广州·“法国姐姐”乔伊丝·乔纳森《小意思》2024巡回演唱会'
$ws.Range("D42").Value = '东风中路299号 广州中山纪念堂'
$ws.Range("E42").Value = '2024.11.17 19:30-11.17 21:00'
$ws.Range("G42").Value = 180
$ws.Range("H42").Value = 'https://show.bilibili.com/platform/detail.html?id=91814'
$ws.Range("I42").Value = '//i2.hdslb.com/bfs/openplatform/202408/bnKPQEEd1725008600562.jpeg'
$ws.Range("F49").Value = 33
